# Adds a "password" column to the "usuario" table, right after "login"
# (new column C), with literal value "None" for every data row, and
# updates the SQL-generating formula in the last column to include the
# new field.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new column before the current column C (primer_apellido).
#    Everything from C onward (including the generator formula column)
#    shifts one column to the right automatically.
$ws.Columns("C:C").Insert()

# 2) Header for the new column.
$ws.Range("C2").Value = "password"

# 3) Literal placeholder value for every existing user row.
$ws.Range("C3:C24").Value = "None"

# 4) Row 3's generator formula is NOT part of the shared formula group,
#    so it is rewritten on its own.
$f3 = @"
="insert into "&B`$1&" ("&CHAR(10)&
CHAR(9)&`$A`$2&","&CHAR(10)&
CHAR(9)&`$B`$2&","&CHAR(10)&
CHAR(9)&`$C`$2&","&CHAR(10)&
CHAR(9)&`$D`$2&","&CHAR(10)&
CHAR(9)&`$E`$2&","&CHAR(10)&
CHAR(9)&`$F`$2&","&CHAR(10)&
CHAR(9)&`$G`$2&","&CHAR(10)&
CHAR(9)&`$H`$2&","&CHAR(10)&
CHAR(9)&`$I`$2&")"&CHAR(10)&
"values ("&CHAR(10)&
CHAR(9)&A3&","&CHAR(10)&
CHAR(9)&"'"&B3&"',"&CHAR(10)&
CHAR(9)&"'"&C3&"',"&CHAR(10)&
CHAR(9)&"'"&D3&"',"&CHAR(10)&
CHAR(9)&"'"&E3&"',"&CHAR(10)&
CHAR(9)&"'"&F3&"',"&CHAR(10)&
CHAR(9)&"'"&G3&"',"&CHAR(10)&
CHAR(9)&H3&","&CHAR(10)&
CHAR(9)&I3&");"&CHAR(10)
"@
$ws.Range("J3").Formula = $f3

# 5) Rows 4:24 share one formula. Writing the same formula text to the
#    whole range at once lets relative references adjust per-row and
#    keeps the cells grouped as a single shared formula, just like the
#    original file.
$f4 = @"
="insert into "&B`$1&" ("&CHAR(10)&
CHAR(9)&`$A`$2&","&CHAR(10)&
CHAR(9)&`$B`$2&","&CHAR(10)&
CHAR(9)&`$C`$2&","&CHAR(10)&
CHAR(9)&`$D`$2&","&CHAR(10)&
CHAR(9)&`$E`$2&","&CHAR(10)&
CHAR(9)&`$F`$2&","&CHAR(10)&
CHAR(9)&`$G`$2&","&CHAR(10)&
CHAR(9)&`$H`$2&","&CHAR(10)&
CHAR(9)&`$I`$2&")"&CHAR(10)&
"values ("&CHAR(10)&
CHAR(9)&A4&","&CHAR(10)&
CHAR(9)&"'"&B4&"',"&CHAR(10)&
CHAR(9)&"'"&C4&"',"&CHAR(10)&
CHAR(9)&"'"&D4&"',"&CHAR(10)&
CHAR(9)&"'"&E4&"',"&CHAR(10)&
CHAR(9)&"'"&F4&"',"&CHAR(10)&
CHAR(9)&"'"&G4&"',"&CHAR(10)&
CHAR(9)&H4&","&CHAR(10)&
CHAR(9)&I4&");"&CHAR(10)
"@
$ws.Range("J4:J24").Formula = $f4

# 6) Row heights grow by one line (password row) for every row except
#    the very first data row.
$ws.Rows("4:24").RowHeight = 236.25

# 7) Give the new column roughly the same width as its neighbours.
$ws.Columns("C:C").ColumnWidth = 36.85546875

# 8) Selection/view cosmetics to match the edited file.
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("J3:J24").Select()
